# atualizado planilha de gastos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fevereiro")

# "IPTU" expense: moved from column C to column B (value unchanged: 631.52)
$ws.Range("C9").Clear()
$ws.Range("B9").Value = 631.52

# "Motorola One" expense: moved from column C to column B, value updated 850 -> 679.9
$ws.Range("C12").Clear()
$ws.Range("B12").Value = 679.9

# New expense row: "PosEAD" = 210.83
$ws.Range("A13").Value = "PosEAD"
$ws.Range("B13").Value = 210.83

# Reflect the new active selection
$ws.Range("B9").Select()
